$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 13 (column S) was removed from the attendance sheet: clear the
# attendance marks (P/A) for every student row (7-82) and restore the
# cell formatting to the sheet's "unused day" look (matching the blank
# column T in the same row), just like clearing out a cancelled class
# column in the sheet.
for ($r = 7; $r -le 82; $r++) {
    $src = $ws.Cells.Item($r, 20)  # column T - already-blank template cell for this row
    $dst = $ws.Cells.Item($r, 19)  # column S - session 13
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

$wb.Save()
